$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-26 Thursday" "2025-06-27 Friday"

Replace-Text "912×5=" "541×6="
Replace-Text "514×3=" "712×3="
Replace-Text "226×6=" "919×9="
Replace-Text "357×7=" "550×9="
Replace-Text "298×6=" "408×3="
Replace-Text "752×5=" "181×7="
Replace-Text "147×4=" "892×3="
Replace-Text "179×6=" "499×6="
Replace-Text "765×4=" "688×2="
Replace-Text "428×3=" "846×7="
Replace-Text "234×2=" "907×5="
Replace-Text "341×6=" "948×2="
Replace-Text "816×8=" "439×7="
Replace-Text "303×7=" "937×3="
Replace-Text "535×4=" "322×4="
Replace-Text "754×9=" "699×9="
Replace-Text "838×3=" "795×4="
Replace-Text "722×5=" "339×7="
Replace-Text "934×5=" "283×8="
Replace-Text "636×7=" "258×9="
Replace-Text "118×6=" "760×9="
Replace-Text "826×2=" "106×4="
Replace-Text "549×4=" "649×6="
Replace-Text "536×9=" "749×8="
Replace-Text "209×6=" "977×2="
